$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Formações")
$ws4.Range("L4").Select() | Out-Null
